$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 1.336873824401267

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 6.048734245549538

$ws.Range("B4").Value = 0.1554434735375247
$ws.Range("C4").Value = 0.05231270169004087
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 3.790242309005201

$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 3.594575437922795

$ws.Range("B6").Value = 0.7287194209349384
$ws.Range("C6").Value = 0.05231270169004087
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 1.993651680714026

$ws.Range("B7").Value = 3.182878228561681
$ws.Range("C7").Value = 9.226618575922256
$ws.Range("D7").Value = 3.082599426703578
$ws.Range("E7").Value = 6.48142807727062
$ws.Range("G7").Value = 21.97352430845813

$ws.Range("B8").Value = 3.182878228561681
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 3.082599426703578
$ws.Range("E8").Value = 0.4998867070740569
$ws.Range("G8").Value = 8.418600821238126

$ws.Range("B9").Value = 3.182878228561681
$ws.Range("C9").Value = 1.65323645889881
$ws.Range("D9").Value = 0.7127328510149897
$ws.Range("E9").Value = 0.4998867070740569
$ws.Range("G9").Value = 6.048734245549538
